$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cryptos list refresh: per-row Price (D) / Volume(1h) (E) updates.
# Rows 40-43 also got reordered (ImmutableX<->Aave, RenderToken<->Kaspa)
# together with their Coin/Link/Price/Volume cells.
#
# Some new Price strings are plain decimals (e.g. "214.94") that Excel's
# usual text -> value coercion would silently turn into numbers, unlike
# the original "1.2.3"-style text prices. Force those specific cells to
# the Text format first so they round-trip as inline/shared strings,
# matching the source data (column D is entirely textual).

$ws.Range('D2').Value = '29.905.34'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').Value = '1.631.62'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('E4').Value = '  +0.67%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.94'
$ws.Range('E5').Value = '  +1.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.522'
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('E7').Value = '  +0.64%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.79'
$ws.Range('E8').Value = '  -1.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.259'
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0899'
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('D12').Value = '1.864.76'
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('D13').Value = '1.637.73'
$ws.Range('E13').Value = '  +1.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.573'
$ws.Range('E14').Value = '  +0.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '9.38'
$ws.Range('E15').Value = '  +4.38%  '
$ws.Range('D16').Value = '29.904.07'
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.84'
$ws.Range('E17').Value = '  -1.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '65.24'
$ws.Range('E18').Value = '  +1.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.88'
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').Value = '0.0₃0702'
$ws.Range('E20').Value = '  -0.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  +0.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.83'
$ws.Range('E22').Value = '  +1.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.14'
$ws.Range('E23').Value = '  +0.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.17'
$ws.Range('E24').Value = '  +2.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.70'
$ws.Range('E25').Value = '  +0.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.50'
$ws.Range('E26').Value = '  -0.75%  '
$ws.Range('E27').Value = '  -0.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.62'
$ws.Range('E28').Value = '  +0.42%  '
$ws.Range('E29').Value = '  +0.54%  '
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('E31').Value = '  +1.55%  '
$ws.Range('E32').Value = '  +2.23%  '
$ws.Range('E33').Value = '  -0.67%  '
$ws.Range('D34').Value = '1.422.75'
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.69'
$ws.Range('E35').Value = '  +3.79%  '
$ws.Range('E36').Value = '  -2.77%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.77'
$ws.Range('E37').Value = '  -3.25%  '
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.558'
$ws.Range('E40').Value = '  +0.39%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '75.17'
$ws.Range('E41').Value = '  +7.81%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.00'
$ws.Range('E42').Value = '  +1.89%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0502'
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.834'
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('E46').Value = '  +0.74%  '
$ws.Range('D47').Value = '1.771.77'
$ws.Range('E47').Value = '  +1.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.33'
$ws.Range('E48').Value = '  -1.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '48.49'
$ws.Range('E49').Value = '  -10.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '91.97'
$ws.Range('E50').Value = '  +4.73%  '
$ws.Range('E51').Value = '  +1.68%  '
